# Adds the "crystal" item/skill entries to the translation sheet.
#
# The original sheet has rows 38-43 holding the "aim_skill_*" and
# "how_to_use_*" entries. The edit inserts two brand new rows (38 and 39)
# with the crystal skill title/description (and pushes everything that used
# to live at rows 38-43 down by 5, to rows 43-48, leaving the same blank-row
# gaps that existed before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 2 new rows while keeping the previous gap pattern
# (old row 40 was blank; after shifting by 5 it should land on row 45).
$ws.Rows("38:42").Insert()

# New row 38: crystal_skill_title / Crystal / Crystal / Crystal
$ws.Range("A38").Value = "crystal_skill_title"

# New row 39: crystal_skill_description / <description> x3
$ws.Range("A39").Value = "crystal_skill_description"

$ws.Range("B38").Value = "Crystal"
$ws.Range("C38").Value = "Crystal"
$ws.Range("D38").Value = "Crystal"

$crystalDescription = "Crystals are not very stable and are fuel for powerfull ancient weapons. Try to put in on fireplace or jsust to throw it. Lets see what will happened. "
$ws.Range("B39").Value = $crystalDescription
$ws.Range("C39").Value = $crystalDescription
$ws.Range("D39").Value = $crystalDescription

# Row heights: match the auto-fit heights Excel produced after the edit.
$ws.Rows(33).RowHeight = 43.2
$ws.Rows(36).RowHeight = 57.6
$ws.Rows(39).RowHeight = 57.6
$ws.Rows(44).RowHeight = 57.6
$ws.Rows(46).RowHeight = 43.2
$ws.Rows(47).RowHeight = 57.6
$ws.Rows(48).RowHeight = 43.2

# Match the final selection/view state from the saved workbook.
[void]$ws.Range("A39").Select()
